$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume table (cryptos list update).
# D-column price cells are written with a leading apostrophe to force
# text interpretation (prevents Excel turning "9.20" into the number 9.2,
# or "43.166.95" into a numeric/date guess), then the style is reset back
# to Normal so no quote-prefix formatting is left behind on the cell.

# Row 2
$ws.Range("D2").Value = "'43.166.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.08%  "

# Row 3
$ws.Range("D3").Value = "'2.276.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.64%  "

# Row 4
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").Value = "'111.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.46%  "

# Row 6
$ws.Range("D6").Value = "'263.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.84%  "

# Row 7
$ws.Range("D7").Value = "'0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.90%  "

# Row 8
$ws.Range("E8").Value = "  +0.27%  "

# Row 9
$ws.Range("D9").Value = "'0.605"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.57%  "

# Row 10
$ws.Range("D10").Value = "'46.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.96%  "

# Row 11
$ws.Range("E11").Value = "  -0.37%  "

# Row 12
$ws.Range("D12").Value = "'9.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.37%  "

# Row 13
$ws.Range("E13").Value = "  +2.02%  "

# Row 14
$ws.Range("D14").Value = "'15.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.51%  "

# Row 15
$ws.Range("D15").Value = "'2.621.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.52%  "

# Row 16
$ws.Range("D16").Value = "'0.861"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.70%  "

# Row 17
$ws.Range("D17").Value = "'2.289.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.42%  "

# Row 18
$ws.Range("D18").Value = "'43.188.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.04%  "

# Row 19
$ws.Range("E19").Value = "  -1.78%  "

# Row 20
$ws.Range("D20").Value = "'6.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.56%  "

# Row 21
$ws.Range("D21").Value = "'71.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.14%  "

# Row 22
$ws.Range("D22").Value = "'2.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.90%  "

# Row 23
$ws.Range("D23").Value = "'233.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.38%  "

# Row 24
$ws.Range("D24").Value = "'9.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.89%  "

# Row 25
$ws.Range("D25").Value = "'2.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.71%  "

# Row 26
$ws.Range("E26").Value = "  +1.86%  "

# Row 27
$ws.Range("D27").Value = "'11.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.40%  "

# Row 28
$ws.Range("E28").Value = "  +0.13%  "

# Row 29
$ws.Range("D29").Value = "'40.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.10%  "

# Row 30
$ws.Range("E30").Value = "  -1.75%  "

# Row 31
$ws.Range("E31").Value = "  -0.85%  "

# Row 32
$ws.Range("D32").Value = "'172.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.17%  "

# Row 33
$ws.Range("D33").Value = "'21.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.60%  "

# Row 34
$ws.Range("D34").Value = "'0.0900"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.34%  "

# Row 35
$ws.Range("D35").Value = "'5.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.18%  "

# Row 36
$ws.Range("E36").Value = "  +1.04%  "

# Row 37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'4.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.84%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0366"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.42%  "

# Row 39
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'4.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.48%  "

# Row 40
$ws.Range("E40").Value = "  -5.37%  "

# Row 41
$ws.Range("E41").Value = "  +6.58%  "

# Row 42
$ws.Range("D42").Value = "'76.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.43%  "

# Row 43
$ws.Range("D43").Value = "'13.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.21%  "

# Row 44
$ws.Range("E44").Value = "  -2.43%  "

# Row 45
$ws.Range("D45").Value = "'6.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.58%  "

# Row 46
$ws.Range("E46").Value = "  +0.17%  "

# Row 47
$ws.Range("E47").Value = "  -4.25%  "

# Row 48
$ws.Range("D48").Value = "'103.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.32%  "

# Row 49
$ws.Range("D49").Value = "'8.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.86%  "

# Row 50
$ws.Range("D50").Value = "'1.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.01%  "

# Row 51
$ws.Range("D51").Value = "'0.0996"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.91%  "
